$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Column by column: header then data value, matching shared-string insert order
$ws.Range("E1").Value = "region"
$ws.Range("E2").Value = "Region-1"

$ws.Range("F1").Value = "fte"
$ws.Range("F2").Value = 1

$ws.Range("G1").Value = "temporaryDepartment"
$ws.Range("G2").Value = "Sub unit-3"

$ws.Range("H1").Value = "bloodGroup"
$ws.Range("H2").Value = "A"

$ws.Range("I1").Value = "hobbies"
$ws.Range("I2").Value = "score goal against nacional"

# Column widths
$ws.Range("E1").EntireColumn.ColumnWidth = 8.85546875
$ws.Range("G1").EntireColumn.ColumnWidth = 21.28515625
$ws.Range("H1").EntireColumn.ColumnWidth = 11.5703125

# Selection
$ws.Range("I2").Select()
